$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "25.268.89"
Set-TextValue $ws "E2" "  -2.94%  "
Set-TextValue $ws "D3" "1.554.78"
Set-TextValue $ws "E3" "  -4.61%  "
Set-TextValue $ws "E4" "  -0.08%  "
Set-TextValue $ws "D5" "207.13"
Set-TextValue $ws "E5" "  -3.31%  "
Set-TextValue $ws "E6" "  -0.09%  "
Set-TextValue $ws "E7" "  -4.99%  "
Set-TextValue $ws "D8" "0.0609"
Set-TextValue $ws "E8" "  -1.76%  "
Set-TextValue $ws "D9" "0.242"
Set-TextValue $ws "E9" "  -3.33%  "
Set-TextValue $ws "D10" "17.74"
Set-TextValue $ws "E10" "  -4.13%  "
Set-TextValue $ws "E11" "  -1.02%  "
Set-TextValue $ws "D12" "1.770.94"
Set-TextValue $ws "E12" "  -4.63%  "
Set-TextValue $ws "B13" "WrappedEther"
Set-TextValue $ws "C13" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws "D13" "1.552.79"
Set-TextValue $ws "E13" "  -4.65%  "
Set-TextValue $ws "B14" "Polkadot"
Set-TextValue $ws "C14" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws "D14" "4.00"
Set-TextValue $ws "E14" "  -4.49%  "
Set-TextValue $ws "E15" "  -4.20%  "
Set-TextValue $ws "D16" "25.264.98"
Set-TextValue $ws "E16" "  -2.98%  "
Set-TextValue $ws "B17" "Litecoin"
Set-TextValue $ws "C17" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws "D17" "58.80"
Set-TextValue $ws "E17" "  -4.46%  "
Set-TextValue $ws "B18" "ShibaInu"
Set-TextValue $ws "C18" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws "D18" "0.0₃0706"
Set-TextValue $ws "E18" "  -4.86%  "
Set-TextValue $ws "E19" "  -0.05%  "
Set-TextValue $ws "D20" "185.52"
Set-TextValue $ws "E20" "  -3.70%  "
Set-TextValue $ws "D21" "4.10"
Set-TextValue $ws "E21" "  -3.56%  "
Set-TextValue $ws "E22" "  -2.92%  "
Set-TextValue $ws "D23" "5.84"
Set-TextValue $ws "E23" "  -3.74%  "
Set-TextValue $ws "E24" "  -3.94%  "
Set-TextValue $ws "D25" "1.01"
Set-TextValue $ws "E25" "  -0.07%  "
Set-TextValue $ws "D26" "140.45"
Set-TextValue $ws "E26" "  -2.69%  "
Set-TextValue $ws "D27" "1.64"
Set-TextValue $ws "E27" "  -4.82%  "
Set-TextValue $ws "D28" "14.85"
Set-TextValue $ws "E28" "  -2.70%  "
Set-TextValue $ws "D29" "6.39"
Set-TextValue $ws "E29" "  -5.04%  "
Set-TextValue $ws "E30" "  -6.90%  "
Set-TextValue $ws "E31" "  -3.35%  "
Set-TextValue $ws "E32" "  -3.30%  "
Set-TextValue $ws "E33" "  -5.06%  "
Set-TextValue $ws "E34" "  -3.21%  "
Set-TextValue $ws "E35" "  -3.77%  "
Set-TextValue $ws "D36" "1.084.85"
Set-TextValue $ws "E36" "  -3.77%  "
Set-TextValue $ws "E37" "  -0.14%  "
Set-TextValue $ws "E38" "  -3.09%  "
Set-TextValue $ws "E39" "  -4.67%  "
Set-TextValue $ws "D40" "0.765"
Set-TextValue $ws "E40" "  -10.12%  "
Set-TextValue $ws "E41" "  -7.94%  "
Set-TextValue $ws "D42" "0.796"
Set-TextValue $ws "E42" "  +5.21%  "
Set-TextValue $ws "D43" "92.45"
Set-TextValue $ws "E43" "  -5.97%  "
Set-TextValue $ws "E44" "  -1.50%  "
Set-TextValue $ws "D45" "1.685.37"
Set-TextValue $ws "E45" "  -4.61%  "
Set-TextValue $ws "E46" "  -2.51%  "
Set-TextValue $ws "E47" "  -1.56%  "
Set-TextValue $ws "D48" "52.32"
Set-TextValue $ws "E48" "  -3.98%  "
Set-TextValue $ws "E49" "  -4.00%  "
Set-TextValue $ws "E50" "  -0.23%  "
Set-TextValue $ws "E51" "  -2.17%  "
